# Added Zeiss CP.3 lenses
# - Appends 14 new Zeiss CP.3 compact prime lens rows to the LensTable sheet
# - Un-hides / re-sizes the handful of rows that the AutoFilter was hiding in
#   the 300-309 band (side effect of the table growing / filter recompute)
# - Grows Table1 to cover the new rows
# - Appends a matching changelog entry

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LensTable")

# ---------------------------------------------------------------------------
# 1. New lens rows (404-417)
# ---------------------------------------------------------------------------
$newLenses = @(
  @{A="Zeiss"; B="CP.3 XD 135mm T2.1 Compact Prime";   C=2017; D=7490; E=1150; F=126.5; G=95; H=2;   I="135"; J="0"; K=0; L=1; M="Full-Frame"}
  @{A="Zeiss"; B="CP.3 25mm T2.1 Compact Prime";        C=2017; D=4390; E=820;  F=83.7;  G=95; H=2;   I="25";  J="0"; K=0; L=1; M="Full-Frame"}
  @{A="Zeiss"; B="CP.3 35mm T2.1 Compact Prime";        C=2017; D=4390; E=800;  F=83.7;  G=95; H=2;   I="35";  J="0"; K=0; L=1; M="Full-Frame"}
  @{A="Zeiss"; B="CP.3 50mm T2.1 Compact Prime";        C=2017; D=4390; E=770;  F=83.7;  G=95; H=2;   I="50";  J="0"; K=0; L=1; M="Full-Frame"}
  @{A="Zeiss"; B="CP.3 85mm T2.1 Compact Prime";        C=2017; D=4390; E=880;  F=83.7;  G=95; H=2;   I="85";  J="0"; K=0; L=1; M="Full-Frame"}
  @{A="Zeiss"; B="CP.3 XD 100mm T2.1 Compact Prime";    C=2017; D=6690; E=1010; F=126.5; G=95; H=2;   I="100"; J="0"; K=0; L=1; M="Full-Frame"}
  @{A="Zeiss"; B="CP.3 XD 15mm T2.9 Compact Prime";     C=2017; D=7490; E=870;  F=83.7;  G=95; H=2.8; I="15";  J="0"; K=0; L=1; M="Full-Frame"}
  @{A="Zeiss"; B="CP.3 XD 18mm T2.9 Compact Prime";     C=2017; D=6690; E=860;  F=83.7;  G=95; H=2.8; I="18";  J="0"; K=0; L=1; M="Full-Frame"}
  @{A="Zeiss"; B="CP.3 XD 21mm T2.9 Compact Prime";     C=2017; D=5790; E=820;  F=83.7;  G=95; H=2.8; I="21";  J="0"; K=0; L=1; M="Full-Frame"}
  @{A="Zeiss"; B="CP.3 XD 25mm T2.1 Compact Prime";     C=2017; D=5790; E=820;  F=83.7;  G=95; H=2;   I="25";  J="0"; K=0; L=1; M="Full-Frame"}
  @{A="Zeiss"; B="CP.3 XD 28mm T2.1 Compact Prime";     C=2017; D=5790; E=840;  F=83.7;  G=95; H=2;   I="28";  J="0"; K=0; L=1; M="Full-Frame"}
  @{A="Zeiss"; B="CP.3 XD 35mm T2.1 Compact Prime";     C=2017; D=5790; E=800;  F=83.7;  G=95; H=2;   I="35";  J="0"; K=0; L=1; M="Full-Frame"}
  @{A="Zeiss"; B="CP.3 XD 50mm T2.1 Compact Prime";     C=2017; D=5790; E=770;  F=83.7;  G=95; H=2;   I="50";  J="0"; K=0; L=1; M="Full-Frame"}
  @{A="Zeiss"; B="CP.3 XD 85mm T2.1 Compact Prime";     C=2017; D=5790; E=880;  F=83.7;  G=95; H=2;   I="85";  J="0"; K=0; L=1; M="Full-Frame"}
)

$startRow = 404
$r = $startRow
foreach ($lens in $newLenses) {
    $ws.Range("A$r").Value = $lens.A
    $ws.Range("B$r").Value = $lens.B
    $ws.Range("C$r").Value = $lens.C
    $ws.Range("D$r").Value = $lens.D
    $ws.Range("E$r").Value = $lens.E
    $ws.Range("F$r").Value = $lens.F
    $ws.Range("G$r").Value = $lens.G
    $ws.Range("H$r").Value = $lens.H
    $ws.Range("I$r").Value = $lens.I
    $ws.Range("J$r").Value = $lens.J
    $ws.Range("K$r").Value = $lens.K
    $ws.Range("L$r").Value = $lens.L
    $ws.Range("M$r").Value = $lens.M
    # "Index" column - continues the running 0-based index sequence
    $ws.Range("Q$r").Value = ($r - 1)
    $r = $r + 1
}
$endRow = $r - 1

# Row 409 (CP.3 XD 100mm) has a taller custom row height in the source file.
$ws.Rows.Item(409).RowHeight = 22.2

# ---------------------------------------------------------------------------
# 2. Row height / visibility tweaks for rows 301, 303-309 (AutoFilter
#    recompute side effect of the new rows being added to the table)
# ---------------------------------------------------------------------------
$rowHeights = [ordered]@{301=13.5; 303=24.9; 304=13.5; 305=18.9; 306=15.9; 307=18.6; 308=12.6; 309=16.2}
foreach ($rr in $rowHeights.Keys) {
    $ws.Rows.Item($rr).RowHeight = $rowHeights[$rr]
    $ws.Rows.Item($rr).Hidden = $false
}

# ---------------------------------------------------------------------------
# 3. Grow Table1 / AutoFilter range to include the new rows
# ---------------------------------------------------------------------------
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:Q$endRow"))

# ---------------------------------------------------------------------------
# 4. Changelog entry
# ---------------------------------------------------------------------------
$cl = $wb.Worksheets.Item("Changelog")

# Copy the date-formatted style down from the row above (so we reuse the
# existing built-in date number format instead of minting a new custom one),
# then set the actual value.
$cl.Range("A47").Copy()
$cl.Range("A48").PasteSpecial(-4122)
$cl.Range("A48").Value = 42910
$cl.Range("B48").Value = "Added Zeiss CP.3 lenses"

$cl.Select()
$cl.Range("G48").Select()

# ---------------------------------------------------------------------------
# 5. View state - scroll / selection to match the saved workbook state.
# LensTable is selected last so it remains the workbook's active tab (as in
# the target file).
# ---------------------------------------------------------------------------
$ws.Select()
$ws.Range("A398").Select()
$excel.ActiveWindow.ScrollRow = 398
$ws.Range("B419").Select()
